$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.647.74"
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").Value = "3.801.45"
$ws.Range("E3").Value = "  +0.86%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.48%  "

$ws.Range("D7").Value = "3.799.71"
$ws.Range("E7").Value = "  +0.79%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  -1.05%  "

$ws.Range("E10").Value = "  -1.30%  "

$ws.Range("E11").Value = "  +2.38%  "

$ws.Range("E12").Value = "  -1.20%  "

$ws.Range("E13").Value = "  -3.51%  "

$ws.Range("E14").Value = "  -2.49%  "

$ws.Range("D15").Value = "4.431.87"
$ws.Range("E15").Value = "  +1.02%  "

$ws.Range("D16").Value = "3.798.17"
$ws.Range("E16").Value = "  +0.84%  "

$ws.Range("D17").Value = "69.672.98"
$ws.Range("E17").Value = "  -0.66%  "

$ws.Range("E18").Value = "  -1.39%  "

$ws.Range("E19").Value = "  -3.98%  "

$ws.Range("E20").Value = "  -1.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "506.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("E22").Value = "  +1.28%  "

$ws.Range("E23").Value = "  +0.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.09%  "

$ws.Range("E26").Value = "  +4.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.59%  "

$ws.Range("E28").Value = "  -5.56%  "

$ws.Range("E30").Value = "  +0.70%  "

$ws.Range("E31").Value = "  +0.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.30%  "

$ws.Range("E34").Value = "  -2.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.04%  "

$ws.Range("E37").Value = "  -2.00%  "

$ws.Range("E38").Value = "  +6.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "481.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.339"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.00%  "

$ws.Range("E41").Value = "  +6.55%  "

$ws.Range("E42").Value = "  -2.69%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "49.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.45%  "

$ws.Range("E45").Value = "  -2.09%  "

$ws.Range("D46").Value = "2.929.81"
$ws.Range("E46").Value = "  -2.77%  "

$ws.Range("E47").Value = "  -1.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "139.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.55%  "

$ws.Range("E49").Value = "  +0.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.81%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.64%  "
